$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 13:46"

# Update province rows whose data changed (values refreshed + table re-sorted
# descending by "Casos totales", which reshuffles which label lands on which row)
$ws.Range("A7").Value = "Bizkaia/Vizcaya"
$ws.Range("B7").Value = 1501
$ws.Range("C7").Value = 466
$ws.Range("D7").Value = 1447
$ws.Range("E7").Value = 54

$ws.Range("A9").Value = "Araba/Alava"
$ws.Range("B9").Value = 1207
$ws.Range("C9").Value = 466
$ws.Range("D9").Value = 1126
$ws.Range("E9").Value = 81

$ws.Range("A12").Value = "Aragon"
$ws.Range("B12").Value = 907
$ws.Range("C12").Value = 29
$ws.Range("D12").Value = 838
$ws.Range("E12").Value = 40

$ws.Range("A19").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B19").Value = 563
$ws.Range("C19").Value = 466
$ws.Range("D19").Value = 543
$ws.Range("E19").Value = 20

$ws.Range("A27").Value = "Albacete"
$ws.Range("B27").Value = 430
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 390
$ws.Range("E27").Value = 32

$ws.Range("A28").Value = "Valladolid"
$ws.Range("B28").Value = 410
$ws.Range("C28").Value = 24
$ws.Range("D28").Value = 369
$ws.Range("E28").Value = 17

$ws.Range("A29").Value = "Tenerife"
$ws.Range("B29").Value = 409
$ws.Range("C29").Value = 15
$ws.Range("D29").Value = 400
$ws.Range("E29").Value = 21

$ws.Range("A44").Value = "Gran Canaria"
$ws.Range("B44").Value = 171
$ws.Range("C44").Value = 15
$ws.Range("D44").Value = 169
$ws.Range("E44").Value = 21

$ws.Range("A53").Value = "Melilla"
$ws.Range("B53").Value = 38
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 37
$ws.Range("E53").Value = 1

$ws.Range("A55").Value = "La Palma"
$ws.Range("B55").Value = 33
$ws.Range("C55").Value = 15
$ws.Range("D55").Value = 33
$ws.Range("E55").Value = 21

$ws.Range("A57").Value = "Fuerteventura"
$ws.Range("B57").Value = 20
$ws.Range("C57").Value = 15
$ws.Range("D57").Value = 20
$ws.Range("E57").Value = 21

$ws.Range("A58").Value = "Lanzarote"
$ws.Range("B58").Value = 17
$ws.Range("C58").Value = 18
$ws.Range("D58").Value = 17
$ws.Range("E58").Value = 21

$ws.Range("A59").Value = "Menorca"
$ws.Range("B59").Value = 15
$ws.Range("C59").Value = 18
$ws.Range("D59").Value = 13
$ws.Range("E59").Value = 0

$ws.Range("A62").Value = "La Gomera"
$ws.Range("B62").Value = 4
$ws.Range("C62").Value = 15
$ws.Range("D62").Value = 2
$ws.Range("E62").Value = 21

$ws.Range("A63").Value = "El Hierro"
$ws.Range("B63").Value = 3
$ws.Range("C63").Value = 15
$ws.Range("D63").Value = 3
$ws.Range("E63").Value = 21
